# [Kadastro App] Yeni kayit eklendi: 2901
# Adds a new record row to both the master "Kayitlar" sheet and the
# corresponding district sheet ("Erdemli"), mirroring the same 6 columns:
# Kayit No, Tarih, Birim, Parsel Sayisi, Is, Personeller.

$wb = $excel.ActiveWorkbook

$kayitNo   = "2901"
$tarih     = "2025-09-08"
$birim     = "Erdemli"
$parsel    = "2"
$is        = "LİHKAB"
$personel  = "ÖZKAN AKBAŞ (Mühendis), SERDAR ARSLAN (Tekniker)"

$sheetNames = @("Kayitlar", "Erdemli")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # New data goes in row 2 (row 1 is the header).
    $row = $ws.Range("A2:F2")

    # Format as text first so values that look numeric/date-like ("2901",
    # "2025-09-08", "2") are stored as literal text, matching the header
    # row's "numberStoredAsText" convention instead of being coerced into
    # numbers / dates.
    $row.NumberFormat = "@"

    $ws.Range("A2").Value = $kayitNo
    $ws.Range("B2").Value = $tarih
    $ws.Range("C2").Value = $birim
    $ws.Range("D2").Value = $parsel
    $ws.Range("E2").Value = $is
    $ws.Range("F2").Value = $personel
}
